$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.559.90"
Set-TextValue $ws.Range("E2") "  +2.18%  "
Set-TextValue $ws.Range("D3") "3.817.58"
Set-TextValue $ws.Range("E3") "  +1.16%  "
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "682.66"
Set-TextValue $ws.Range("E5") "  +8.68%  "
Set-TextValue $ws.Range("D6") "170.88"
Set-TextValue $ws.Range("E6") "  +3.27%  "
Set-TextValue $ws.Range("D7") "3.815.36"
Set-TextValue $ws.Range("E7") "  +1.10%  "
Set-TextValue $ws.Range("E8") "  +0.03%  "
Set-TextValue $ws.Range("D9") "0.526"
Set-TextValue $ws.Range("E9") "  +0.88%  "
Set-TextValue $ws.Range("E10") "  +1.93%  "
Set-TextValue $ws.Range("D11") "7.23"
Set-TextValue $ws.Range("E11") "  +6.96%  "
Set-TextValue $ws.Range("E12") "  +0.68%  "
Set-TextValue $ws.Range("E13") "  +0.03%  "
Set-TextValue $ws.Range("D14") "35.94"
Set-TextValue $ws.Range("E14") "  +1.65%  "
Set-TextValue $ws.Range("D15") "4.459.03"
Set-TextValue $ws.Range("E15") "  +1.15%  "
Set-TextValue $ws.Range("D16") "3.816.87"
Set-TextValue $ws.Range("E16") "  +1.67%  "
Set-TextValue $ws.Range("D17") "70.591.21"
Set-TextValue $ws.Range("E17") "  +2.16%  "
Set-TextValue $ws.Range("D18") "17.72"
Set-TextValue $ws.Range("E18") "  +0.76%  "
Set-TextValue $ws.Range("E19") "  +2.27%  "
Set-TextValue $ws.Range("E20") "  +0.62%  "
Set-TextValue $ws.Range("D21") "11.13"
Set-TextValue $ws.Range("E21") "  +16.66%  "
Set-TextValue $ws.Range("D22") "476.83"
Set-TextValue $ws.Range("E22") "  +1.99%  "
Set-TextValue $ws.Range("D23") "0.715"
Set-TextValue $ws.Range("E23") "  +1.29%  "
Set-TextValue $ws.Range("D24") "83.42"
Set-TextValue $ws.Range("E24") "  +0.53%  "
Set-TextValue $ws.Range("E25") "  -0.94%  "
Set-TextValue $ws.Range("E26") "  +2.13%  "
Set-TextValue $ws.Range("D28") "2.12"
Set-TextValue $ws.Range("E28") "  -1.73%  "
Set-TextValue $ws.Range("E29") "  +0.05%  "
Set-TextValue $ws.Range("D30") "3.967.54"
Set-TextValue $ws.Range("E30") "  +1.13%  "
Set-TextValue $ws.Range("E31") "  +9.82%  "
Set-TextValue $ws.Range("E32") "  +2.73%  "
Set-TextValue $ws.Range("E33") "  +4.01%  "
Set-TextValue $ws.Range("D34") "29.60"
Set-TextValue $ws.Range("E34") "  +2.95%  "
Set-TextValue $ws.Range("D35") "0.181"
Set-TextValue $ws.Range("E35") "  +4.25%  "
Set-TextValue $ws.Range("E36") "  +2.22%  "
Set-TextValue $ws.Range("E37") "  +0.07%  "
Set-TextValue $ws.Range("D38") "3.772.52"
Set-TextValue $ws.Range("E38") "  +1.28%  "
Set-TextValue $ws.Range("E39") "  +1.12%  "
Set-TextValue $ws.Range("D40") "3.40"
Set-TextValue $ws.Range("E40") "  +2.68%  "
Set-TextValue $ws.Range("D41") "5.95"
Set-TextValue $ws.Range("E41") "  +2.26%  "
Set-TextValue $ws.Range("D42") "0.965"
Set-TextValue $ws.Range("E42") "  -0.34%  "
Set-TextValue $ws.Range("E43") "  -0.06%  "
Set-TextValue $ws.Range("E44") "  +11.64%  "
Set-TextValue $ws.Range("D46") "46.02"
Set-TextValue $ws.Range("E46") "  +6.71%  "
Set-TextValue $ws.Range("D47") "159.41"
Set-TextValue $ws.Range("E47") "  +3.03%  "
Set-TextValue $ws.Range("E48") "  +11.10%  "
Set-TextValue $ws.Range("D51") "0.301"
Set-TextValue $ws.Range("E51") "  +1.91%  "

# Row 49/50: OKB and ONDO swapped positions
Set-TextValue $ws.Range("B49") "ONDO"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D49") "1.45"
Set-TextValue $ws.Range("E49") "  +6.89%  "
Set-TextValue $ws.Range("B50") "OKB"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D50") "48.19"
Set-TextValue $ws.Range("E50") "  +3.17%  "
